# ============================================================
# Edit: insert "2022-Q4" quarter
#  1. "总计" (summary) sheet: insert a new row 2 with the 2022-Q4 totals,
#     pushing all existing quarters down by one row.
#  2. Insert a brand-new worksheet named "2022-Q4" right after "总计"
#     (i.e. before the current first quarter tab) containing the fund
#     holdings detail for the new quarter.
# ============================================================

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) "总计" sheet - prepend a row for 2022-Q4
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Rows.Item(2).Insert()

# Copy the row-2 (now row-3, the old "2022-Q3" row) formatting up into
# the freshly inserted blank row 2 so the new row matches the sheet's
# existing look (bold/centered/bordered "A" column, plain data cells).
$summary.Range("A3:D3").Copy()
$summary.Range("A2:D2").PasteSpecial(-4122)

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 21
$summary.Cells.Item(2, 4).Value = 2.43

# ------------------------------------------------------------------
# 2) New "2022-Q4" worksheet with the fund holdings detail
# ------------------------------------------------------------------
$q4data = @(
    @(0, "009010", "华夏兴阳一年持有期混合", "27.18", "90.65", "3.32", "0.9024", 5),
    @(1, "012434", "银华多元回报一年持有期混合", "20.31", "87.60", "1.79", "0.3635", 6),
    @(2, "001678", "英大国企改革主题股票", "4.20", "92.20", "8.11", "0.3406", 2),
    @(3, "003713", "英大睿盛灵活配置混合A", "2.39", "93.29", "8.63", "0.2063", 5),
    @(4, "003714", "英大睿盛灵活配置混合C", "2.39", "93.29", "8.63", "0.2063", 5),
    @(5, "003655", "信澳新财富灵活配置混合", "3.03", "78.65", "2.07", "0.0627", 8),
    @(6, "377150", "上投摩根健康品质生活混合A", "3.25", "79.00", "1.90", "0.0618", 10),
    @(7, "015346", "上投摩根健康品质生活混合C", "2.77", "79.00", "1.90", "0.0526", 10),
    @(8, "001607", "英大策略优选混合A", "0.59", "93.12", "8.11", "0.0478", 2),
    @(9, "180028", "银华永祥灵活配置混合", "0.70", "77.51", "5.24", "0.0367", 3),
    @(10, "002307", "银华多元视野灵活配置混合", "1.52", "89.13", "2.10", "0.0319", 9),
    @(11, "003446", "英大睿鑫灵活配置混合A", "0.27", "93.18", "8.23", "0.0222", 3),
    @(12, "012522", "英大稳固增强核心一年持有混合C", "1.05", "23.17", "2.06", "0.0216", 3),
    @(13, "003447", "英大睿鑫灵活配置混合C", "0.22", "93.18", "8.23", "0.0181", 3),
    @(14, "012005", "信澳恒盛混合A", "1.29", "35.97", "1.11", "0.0143", 4),
    @(15, "012521", "英大稳固增强核心一年持有混合A", "0.63", "23.17", "2.06", "0.0130", 3),
    @(16, "012854", "英大中证ESG120策略指数A", "0.38", "93.90", "2.94", "0.0112", 5),
    @(17, "002005", "工银新得利混合", "0.49", "26.23", "2.17", "0.0106", 3),
    @(18, "001608", "英大策略优选混合C", "0.02", "93.12", "8.11", "0.0016", 2),
    @(19, "012006", "信澳恒盛混合C", "0.09", "35.97", "1.11", "0.0010", 4),
    @(20, "012855", "英大中证ESG120策略指数C", "0.01", "93.90", "2.94", "0.0003", 5)
)


$firstQuarter = $wb.Worksheets.Item(2)
$q4 = $wb.Worksheets.Add($firstQuarter)
$q4.Name = "2022-Q4"

# Borrow the header + data-row formatting from the existing "2022-Q3"
# sheet (same column layout: A index | B..G fund facts | H rank).
$template = $wb.Worksheets.Item("2022-Q3")

$template.Range("A1:H1").Copy()
$q4.Range("A1:H1").PasteSpecial(-4122)

$template.Range("A2:H2").Copy()
$q4.Range("A2:H22").PasteSpecial(-4122)

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 0; $col -lt $headers.Length; $col++) {
    $q4.Cells.Item(1, 2 + $col).Value = $headers[$col]
}

for ($i = 0; $i -lt $q4data.Length; $i++) {
    $item = $q4data[$i]
    $r = 2 + $i

    $q4.Cells.Item($r, 1).Value = $item[0]

    $textRange = $q4.Range($q4.Cells.Item($r, 2), $q4.Cells.Item($r, 7))
    $q4.Cells.Item($r, 2).Value = "'" + $item[1]
    $q4.Cells.Item($r, 3).Value = "'" + $item[2]
    $q4.Cells.Item($r, 4).Value = "'" + $item[3]
    $q4.Cells.Item($r, 5).Value = "'" + $item[4]
    $q4.Cells.Item($r, 6).Value = "'" + $item[5]
    $q4.Cells.Item($r, 7).Value = "'" + $item[6]
    $textRange.Style = "Normal"

    $q4.Cells.Item($r, 8).Value = $item[7]
}

Write-Host "2022-Q4 sheet + summary row inserted."
